# Generate Report for Handoff
#
# This script updates the localization-status report after a new handoff
# was generated for the 798d0941 e2e file:
#   - Overview sheet: zh-cn / de-de status columns flip from
#     "Handed back: in sync with en-US" to "Ready for handoff", and the
#     "Latest HO Xliff Generate Date" timestamp advances.
#   - zh-cn / de-de sheets: Priority flips from "ht" to "mt", the
#     "Latest Handoff Datetime" timestamp advances, and the second data
#     row (798d0941 file) now reports a stale-handback Error Detail
#     pointing at the two git blobs.
#   - A handful of columns are narrowed/widened to fit the new text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"

$oldOverviewDate = "2016-10-21 04:53:42"
$newOverviewDate = "2016-10-21 04:55:30"

$oldPriority = "ht"
$newPriority = "mt"

$oldHandoffDate = "2016-10-21 04:53:10"
$newHandoffDate = "2016-10-21 04:55:18"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/44a5bb7bb5c03a042accccf1db8113adeeb2abdc/e2e/798d0941-2185-41e8-8db3-66ec1f5541b2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5e1379e62979cecfd6439f1fc13b71669d139ce5/e2e/798d0941-2185-41e8-8db3-66ec1f5541b2.md."

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$overview.Range("G2").Value = $newOverviewDate
$overview.Range("G3").Value = $newOverviewDate

# zh-cn / de-de status columns got narrower now that the text is shorter
$overview.Columns.Item(5).ColumnWidth = 16.333333333333336
$overview.Columns.Item(6).ColumnWidth = 16.333333333333336

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("E2").Value = $newPriority
$zhcn.Range("E3").Value = $newPriority

$zhcn.Range("H2").Value = $newHandoffDate
$zhcn.Range("H3").Value = $newHandoffDate

$zhcn.Range("P3").Value = $errorDetail

$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333336
$zhcn.Columns.Item(16).ColumnWidth = 39.16666666666667

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("E2").Value = $newPriority
$dede.Range("E3").Value = $newPriority

$dede.Range("H2").Value = $newOverviewDate
$dede.Range("H3").Value = $newOverviewDate

$dede.Range("P3").Value = $errorDetail

$dede.Columns.Item(3).ColumnWidth = 16.333333333333336
$dede.Columns.Item(16).ColumnWidth = 39.16666666666667
